# Warrant Issued Report - add a new mapping row for the
# "Vehicle Registration Non-Expiring Indicator" element.
#
# The new row is inserted just above the existing "Vehicle RegistrationState"
# row (which currently lives at row 50), so it lands at row 50 and pushes
# everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warrant Issued Report")

$newRow = 50

# Push the existing row 50 ("Vehicle RegistrationState") and everything
# below it down by one row, opening up a blank row 50 for the new entry.
$ws.Rows.Item($newRow).Insert()

# Fill in the new mapping entry.
$ws.Cells.Item($newRow, 3).Value = "Vehicle Registration Non-Expiring Indicator"
$ws.Cells.Item($newRow, 4).Value = "True if vehicle registration is non-expiring"
$ws.Cells.Item($newRow, 5).Value = $true
$ws.Cells.Item($newRow, 6).Value = "/wir-doc:WarrantIssuedReport/j:ConveyanceRegistration[not(j:RegistrationExpirationDate)]/wir-ext:ConveyanceRegistrationNonExpiringIndicator"

# Match the look of the surrounding mapping rows: wrapped text, top-left
# aligned, tall enough to show the full XPath expression.
$dataRange = $ws.Range("C" + $newRow + ":F" + $newRow)
$dataRange.WrapText = $true
$dataRange.VerticalAlignment = -4160
$dataRange.HorizontalAlignment = -4131
$ws.Rows.Item($newRow).RowHeight = 56

# Leave the cursor near the new entry, same general area the author was
# working in when they saved.
$ws.Range("F45").Select() | Out-Null

Write-Host "Inserted row $newRow with the Vehicle Registration Non-Expiring Indicator mapping."
